$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new quarterly columns S (31/03/2024) and T (30/06/2024) ---
# Copy formatting (bold, centered, bordered) from the existing header cell R1
# onto the two new header cells before writing their labels.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("S1").Value = "31/03/2024"
$ws.Range("T1").Value = "30/06/2024"

# --- Data rows 2-80: append the two new quarters figures ---
# Row 2
$ws.Cells.Item(2, 19).Value = 648003.008
$ws.Cells.Item(2, 20).Value = 646036.992
# Row 3
$ws.Cells.Item(3, 19).Value = 257724.992
$ws.Cells.Item(3, 20).Value = 249991.008
# Row 4
$ws.Cells.Item(4, 19).Value = 59644
$ws.Cells.Item(4, 20).Value = 60363
# Row 5
$ws.Cells.Item(5, 19).Value = 109578
$ws.Cells.Item(5, 20).Value = 106386
# Row 6
$ws.Cells.Item(6, 19).Value = 68538
$ws.Cells.Item(6, 20).Value = 67371
# Row 7
$ws.Cells.Item(7, 19).Value = 0
$ws.Cells.Item(7, 20).Value = 0
# Row 8
$ws.Cells.Item(8, 19).Value = 0
$ws.Cells.Item(8, 20).Value = 0
# Row 9
$ws.Cells.Item(9, 19).Value = 15763
$ws.Cells.Item(9, 20).Value = 7660
# Row 10
$ws.Cells.Item(10, 19).Value = 4142
$ws.Cells.Item(10, 20).Value = 8151
# Row 11
$ws.Cells.Item(11, 19).Value = 60
$ws.Cells.Item(11, 20).Value = 60
# Row 12
$ws.Cells.Item(12, 19).Value = 6550
$ws.Cells.Item(12, 20).Value = 14879
# Row 13
$ws.Cells.Item(13, 19).Value = 0
$ws.Cells.Item(13, 20).Value = 0
# Row 14
$ws.Cells.Item(14, 19).Value = 0
$ws.Cells.Item(14, 20).Value = 0
# Row 15
$ws.Cells.Item(15, 19).Value = 0
$ws.Cells.Item(15, 20).Value = 0
# Row 16
$ws.Cells.Item(16, 19).Value = 50
$ws.Cells.Item(16, 20).Value = 168
# Row 17
$ws.Cells.Item(17, 19).Value = 0
$ws.Cells.Item(17, 20).Value = 0
# Row 18
$ws.Cells.Item(18, 19).Value = 0
$ws.Cells.Item(18, 20).Value = 0
# Row 19
$ws.Cells.Item(19, 19).Value = 6323
$ws.Cells.Item(19, 20).Value = 6046
# Row 20
$ws.Cells.Item(20, 19).Value = 177
$ws.Cells.Item(20, 20).Value = 513
# Row 21
$ws.Cells.Item(21, 19).Value = 0
$ws.Cells.Item(21, 20).Value = 0
# Row 22
$ws.Cells.Item(22, 19).Value = 0
$ws.Cells.Item(22, 20).Value = 0
# Row 23
$ws.Cells.Item(23, 19).Value = 14822
$ws.Cells.Item(23, 20).Value = 13912
# Row 24
$ws.Cells.Item(24, 19).Value = 368905.984
$ws.Cells.Item(24, 20).Value = 367255.008
# Row 25
$ws.Cells.Item(25, 19).Value = 0
$ws.Cells.Item(25, 20).Value = 0
# Row 26
$ws.Cells.Item(26, 19).Value = 648003.008
$ws.Cells.Item(26, 20).Value = 646036.992
# Row 27
$ws.Cells.Item(27, 19).Value = 99145
$ws.Cells.Item(27, 20).Value = 104566
# Row 28
$ws.Cells.Item(28, 19).Value = 21016
$ws.Cells.Item(28, 20).Value = 27601
# Row 29
$ws.Cells.Item(29, 19).Value = 22492
$ws.Cells.Item(29, 20).Value = 26877
# Row 30
$ws.Cells.Item(30, 19).Value = 7737
$ws.Cells.Item(30, 20).Value = 6694
# Row 31
$ws.Cells.Item(31, 19).Value = 7503
$ws.Cells.Item(31, 20).Value = 4970
# Row 32
$ws.Cells.Item(32, 19).Value = 0
$ws.Cells.Item(32, 20).Value = 0
# Row 33
$ws.Cells.Item(33, 19).Value = 0
$ws.Cells.Item(33, 20).Value = 0
# Row 34
$ws.Cells.Item(34, 19).Value = 40397
$ws.Cells.Item(34, 20).Value = 38424
# Row 35
$ws.Cells.Item(35, 19).Value = 0
$ws.Cells.Item(35, 20).Value = 0
# Row 36
$ws.Cells.Item(36, 19).Value = 0
$ws.Cells.Item(36, 20).Value = 0
# Row 37
$ws.Cells.Item(37, 19).Value = 100733
$ws.Cells.Item(37, 20).Value = 94645
# Row 38
$ws.Cells.Item(38, 19).Value = 13048
$ws.Cells.Item(38, 20).Value = 12901
# Row 39
$ws.Cells.Item(39, 19).Value = 262
$ws.Cells.Item(39, 20).Value = 0
# Row 40
$ws.Cells.Item(40, 19).Value = 15589
$ws.Cells.Item(40, 20).Value = 14405
# Row 41
$ws.Cells.Item(41, 19).Value = 70184
$ws.Cells.Item(41, 20).Value = 65807
# Row 42
$ws.Cells.Item(42, 19).Value = 0
$ws.Cells.Item(42, 20).Value = 0
# Row 43
$ws.Cells.Item(43, 19).Value = 1650
$ws.Cells.Item(43, 20).Value = 1532
# Row 44
$ws.Cells.Item(44, 19).Value = 0
$ws.Cells.Item(44, 20).Value = 0
# Row 45
$ws.Cells.Item(45, 19).Value = 0
$ws.Cells.Item(45, 20).Value = 0
# Row 46
$ws.Cells.Item(46, 19).Value = 0
$ws.Cells.Item(46, 20).Value = 0
# Row 47
$ws.Cells.Item(47, 19).Value = 448124.992
$ws.Cells.Item(47, 20).Value = 446825.984
# Row 48
$ws.Cells.Item(48, 19).Value = 440108
$ws.Cells.Item(48, 20).Value = 440108
# Row 49
$ws.Cells.Item(49, 19).Value = -41180
$ws.Cells.Item(49, 20).Value = -41180
# Row 50
$ws.Cells.Item(50, 19).Value = 0
$ws.Cells.Item(50, 20).Value = 0
# Row 51
$ws.Cells.Item(51, 19).Value = 51202
$ws.Cells.Item(51, 20).Value = 51096
# Row 52
$ws.Cells.Item(52, 19).Value = -3410
$ws.Cells.Item(52, 20).Value = -9024
# Row 53
$ws.Cells.Item(53, 19).Value = 0
$ws.Cells.Item(53, 20).Value = 0
# Row 54
$ws.Cells.Item(54, 19).Value = 1405
$ws.Cells.Item(54, 20).Value = 5826
# Row 55
$ws.Cells.Item(55, 19).Value = 0
$ws.Cells.Item(55, 20).Value = 0
# Row 56
$ws.Cells.Item(56, 19).Value = 0
$ws.Cells.Item(56, 20).Value = 0
# Row 57
$ws.Cells.Item(57, 19).Value = "'"
$ws.Cells.Item(57, 20).Value = "'"
# Row 58
$ws.Cells.Item(58, 19).Value = "'"
$ws.Cells.Item(58, 20).Value = "'"
# Row 59
$ws.Cells.Item(59, 19).Value = 67980
$ws.Cells.Item(59, 20).Value = 67814
# Row 60
$ws.Cells.Item(60, 19).Value = -27221
$ws.Cells.Item(60, 20).Value = -28043
# Row 61
$ws.Cells.Item(61, 19).Value = 40759
$ws.Cells.Item(61, 20).Value = 39771
# Row 62
$ws.Cells.Item(62, 19).Value = -12520
$ws.Cells.Item(62, 20).Value = -12108
# Row 63
$ws.Cells.Item(63, 19).Value = -31761
$ws.Cells.Item(63, 20).Value = -34835
# Row 64
$ws.Cells.Item(64, 19).Value = 0
$ws.Cells.Item(64, 20).Value = 0
# Row 65
$ws.Cells.Item(65, 19).Value = 0
$ws.Cells.Item(65, 20).Value = 0
# Row 66
$ws.Cells.Item(66, 19).Value = -709
$ws.Cells.Item(66, 20).Value = 182
# Row 67
$ws.Cells.Item(67, 19).Value = 0
$ws.Cells.Item(67, 20).Value = 0
# Row 68
$ws.Cells.Item(68, 19).Value = 1636
$ws.Cells.Item(68, 20).Value = 376
# Row 69
$ws.Cells.Item(69, 19).Value = 4307
$ws.Cells.Item(69, 20).Value = 2556
# Row 70
$ws.Cells.Item(70, 19).Value = -2671
$ws.Cells.Item(70, 20).Value = -2180
# Row 71
$ws.Cells.Item(71, 19).Value = "'"
$ws.Cells.Item(71, 20).Value = "'"
# Row 72
$ws.Cells.Item(72, 19).Value = "'"
$ws.Cells.Item(72, 20).Value = "'"
# Row 73
$ws.Cells.Item(73, 19).Value = "'"
$ws.Cells.Item(73, 20).Value = "'"
# Row 74
$ws.Cells.Item(74, 19).Value = -2595
$ws.Cells.Item(74, 20).Value = -6614
# Row 75
$ws.Cells.Item(75, 19).Value = -958
$ws.Cells.Item(75, 20).Value = -1083
# Row 76
$ws.Cells.Item(76, 19).Value = 143
$ws.Cells.Item(76, 20).Value = 2704
# Row 77
$ws.Cells.Item(77, 19).Value = "'"
$ws.Cells.Item(77, 20).Value = "'"
# Row 78
$ws.Cells.Item(78, 19).Value = "'"
$ws.Cells.Item(78, 20).Value = "'"
# Row 79
$ws.Cells.Item(79, 19).Value = 0
$ws.Cells.Item(79, 20).Value = 0
# Row 80
$ws.Cells.Item(80, 19).Value = -3410
$ws.Cells.Item(80, 20).Value = -4993
